$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 0) Capture the formatting template of the bold "Meta description"
#    run (leading empty run + a run with only <w:b/>). We need it to
#    build the new bold heading paragraph near the end of the
#    document. NOTE: Range.FormattedText is a *live*, position-bound
#    object (not a detached copy), so we must finish using it before
#    we delete/shift the paragraph it points at.
# ------------------------------------------------------------------
$metaFind = $d.Content.Find
$metaFind.ClearFormatting()
$metaFind.Text = "Meta description"
$metaFind.Forward = $true
$metaFind.Wrap = 0
$metaFound = $metaFind.Execute()
$boldTemplate = $null
if ($metaFound) {
    $boldSrcRange = $d.Range($metaFind.Parent.Start, $metaFind.Parent.End)
    $boldTemplate = $boldSrcRange.FormattedText
}

# ------------------------------------------------------------------
# 1) Insert a new bold "Play Free Book of Kings Slot - Big Wins
#    Await" paragraph right before the trailing AI-image-prompt
#    paragraph (anchored after the preceding "Limited paylines ..."
#    bullet -- a plain, non-italic paragraph -- so the new run does
#    not inherit italic formatting). The paragraph style is then
#    dropped back to Normal and the captured bold-run formatting is
#    stamped onto the new text, reproducing the leading empty run +
#    bold run structure.
# ------------------------------------------------------------------
$bulletFind = $d.Content.Find
$bulletFind.ClearFormatting()
$bulletFind.Text = "Limited paylines may not provide enough variety"
$bulletFind.Forward = $true
$bulletFind.Wrap = 0
if ($bulletFind.Execute()) {
    $bulletRange = $d.Range($bulletFind.Parent.Start, $bulletFind.Parent.End)

    $insertAnchor = $bulletRange.Duplicate
    $insertAnchor.Collapse(0)
    $insertAnchor.InsertParagraphAfter()

    $paraCount = $d.Paragraphs.Count
    $newHeadingPara = $d.Paragraphs.Item($paraCount - 1)
    $newHeadingPara.Style = "Normal"

    $insPt = $newHeadingPara.Range.Duplicate
    $insPt.Collapse(1)

    if ($boldTemplate -ne $null) {
        $insPt.FormattedText = $boldTemplate
        $newHeadingPara2 = $d.Paragraphs.Item($paraCount - 1)
        $newHeadingPara2.Range.Find.Execute("Meta description", $false, $false, $false, $false, $false, $true, 1, $false, "Play Free Book of Kings Slot - Big Wins Await", 2) | Out-Null
    } else {
        $insPt.InsertAfter("Play Free Book of Kings Slot - Big Wins Await")
        $headingStart = $newHeadingPara.Range.Start
        $headingEnd = $headingStart + ("Play Free Book of Kings Slot - Big Wins Await").Length
        $d.Range($headingStart, $headingEnd).Bold = 1
    }
}

# ------------------------------------------------------------------
# 2) Turn the old AI-image "Prompt: ..." paragraph's text into the
#    meta-description copy (keeps its existing italic run/
#    formatting in place).
# ------------------------------------------------------------------
$oldPrompt = 'Prompt: Design a feature image for the online slot game "Book of Kings" that captures the excitement and adventure of discovering ancient treasures in Egypt. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The warrior should be holding a book and standing in front of a pyramid while surrounded by symbols of adventure such as a compass, a treasure chest, and perhaps some hieroglyphics. The colors should be warm and inviting, with shades of gold, brown, and orange. The overall vibe should be one of excitement and possibility, enticing players to explore the riches of the game.'
$newDescription = 'Discover the mysteries of Ancient Egypt and play Book of Kings, an online slot game with free bonuses, high rewards, and mobile compatibility. Play free now.'
$d.Content.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newDescription, 2) | Out-Null

# ------------------------------------------------------------------
# 3) Remove the whole original "Meta description: ..." paragraph
#    that followed the title heading.
# ------------------------------------------------------------------
$metaParaFind = $d.Content.Find
$metaParaFind.ClearFormatting()
$metaParaFind.Text = "Meta description"
$metaParaFind.Forward = $true
$metaParaFind.Wrap = 0
if ($metaParaFind.Execute()) {
    $metaRange = $d.Range($metaParaFind.Parent.Start, $metaParaFind.Parent.End)
    $metaRange.Expand(4) | Out-Null
    $metaRange.Delete()
}

Write-Output "done"
